$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholder auto-text: 1/27/2020 -> 1/31/2020
#    (slide master, all slide layouts, and the notes master each carry an
#    independent <a:fld type="datetimeFigureOut"> run with this cached text)
# ---------------------------------------------------------------------------

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "1/27/2020") {
                $tr.Text = "1/31/2020"
            }
        }
    }
}

# Slide master
Update-DateShape $p.SlideMaster.Shapes

# Every custom (slide) layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

# Notes master
Update-DateShape $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 12: correction in description of how the game board is
#    initialized. Replace the explanatory paragraph and drop the
#    (now redundant / struck-through) trailing paragraphs.
# ---------------------------------------------------------------------------

$s12 = $p.Slides.Item(12)
$shp = $s12.Shapes.Item(2)   # "TextBox 4"

$run1 = "You" + [char]0x2019 + "ll want to declare and initialize the values of these 2D arrays after the line " + [char]0x201C + "public class "
$run2 = "TileGameApp"
$run3 = " {" + [char]0x2026 + [char]0x201D + " and before the main method. "

$tr = $shp.TextFrame.TextRange
$tr.Text = $run1 + $run2 + $run3

# Re-split into three runs so "TileGameApp" keeps standing on its own run,
# matching the original run layout (formatting-wise unaffected).
$len1 = $run1.Length
$len2 = $run2.Length
$len3 = $run3.Length

$part1 = $tr.Characters(1, $len1)
$part1.Text = $run1

$part2 = $tr.Characters($len1 + 1, $len2)
$part2.Text = $run2

$part3 = $tr.Characters($len1 + $len2 + 1, $len3)
$part3.Text = $run3
